$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.86069999999999
$ws.Range("C3").Value = -10.93219999999998
$ws.Range("A4").Value = -21.9566
$ws.Range("B4").Value = 5.621699999999996
$ws.Range("C4").Value = -13.64350000000001
$ws.Range("B5").Value = 5.655999999999995
$ws.Range("A6").Value = -21.2931
$ws.Range("A7").Value = -21.22389999999999
$ws.Range("B8").Value = 5.0758
$ws.Range("C9").Value = -11.57380000000001
$ws.Range("C11").Value = -14.11300000000001
$ws.Range("C14").Value = -11.9399
$ws.Range("A16").Value = -21.22230000000001
$ws.Range("B16").Value = 5.230300000000002
$ws.Range("C18").Value = -14.58360000000001
$ws.Range("A20").Value = -22.80190000000001
$ws.Range("E20").Value = 12.4368
$ws.Range("B22").Value = 5.521999999999999
$ws.Range("C25").Value = -10.93859999999999
